$d = $word.ActiveDocument

# --- Change 1: Append " and expand it's selection" (with proofErr markers) to the
#     "Click on each button ... each button should work as expected" paragraph.
$p1 = $d.Paragraphs(11)
Write-Host "Before: $($p1.Range.Text)"
$frag1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Click on each button </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> each button should work as expected</w:t></w:r><w:r><w:t xml:space="preserve"> and expand </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>it&#8217;s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> selection</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p1.Range.InsertXML($frag1)
Write-Host "After: $($d.Paragraphs(11).Range.Text)"

# --- Change 2: Remove the "Expand each section present..." paragraph entirely.
$p2 = $d.Paragraphs(12)
Write-Host "Deleting: $($p2.Range.Text)"
$p2.Range.Delete()

Write-Host "--- final check ---"
for ($i = 9; $i -le 14; $i++) {
    $pp = $d.Paragraphs($i)
    Write-Host "P$i`: $($pp.Range.Text)"
}
